# Fix the "harvester" column (B) in rnaSamples: Holly added S.GISH to
# harvester in bioSamples, so propagate that value down the harvester
# column here as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B22").Value = "S.GISH"

# Leave the selection on the harvester column, matching where the edit
# was made.
$ws.Columns("B:B").Select()
